$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "0.3.0"
$ws.Range("B16").Value = "2021-04-07T00:00:00+00:00"
$ws.Range("E132").Value = "http://purl.obolibrary.org/obo/OGMS_0000090`n"
$ws.Range("H132").Value = "http://purl.obolibrary.org/obo/OGMS_0000090`n"
$ws.Range("H137").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&version=21.03e&ns=ncit&code=C16809&key=1826977516&b=1&n=null"
$ws.Range("H138").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&version=21.03e&ns=ncit&code=C38101&key=1939144925&b=1&n=null"
$ws.Range("H139").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&version=21.03e&ns=ncit&code=C64384&key=168018700&b=1&n=null"
$ws.Range("H140").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&version=21.03e&ns=ncit&code=C17204&key=1418861988&b=1&n=null"
$ws.Range("H141").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C62663"
$ws.Range("H148").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&version=21.03e&ns=ncit&code=C43433&key=n1963439227&b=1&n=null"
$ws.Range("H151").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&version=21.03e&ns=ncit&code=C54273&key=1753678631&b=1&n=null"
$ws.Range("E152").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C91063"
$ws.Range("H152").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C91063"
$ws.Range("E165").Value = "http://purl.bioontology.org/ontology/MESH/D008907"
$ws.Range("E166").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14238"
$ws.Range("H166").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14238"
$ws.Range("E167").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C160998"
$ws.Range("H167").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C160998"
$ws.Range("E168").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14283"
$ws.Range("H168").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14283"
$ws.Range("E169").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14187"
$ws.Range("H169").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14187"
$ws.Range("E170").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C28176"
$ws.Range("E171").Value = "https://ncit.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&version=21.02d&ns=ncit&code=C45306&key=n1938058024&b=1&n=null"
$ws.Range("E172").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C12434"
$ws.Range("H172").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C12434"
$ws.Range("E173").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C94552"
$ws.Range("H173").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C94552"
$ws.Range("E174").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13234"
$ws.Range("E175").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13356"
$ws.Range("H175").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13356,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/50863008"
$ws.Range("E176").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13275"
$ws.Range("H176").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13275,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/256897009 "
$ws.Range("E177").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13277"
$ws.Range("H177").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13277"
$ws.Range("E178").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13325"
$ws.Range("H178").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13325,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/67922002"
$ws.Range("E179").Value = "http://purl.obolibrary.org/obo/ENVO_00002018"
$ws.Range("H179").Value = "http://purl.obolibrary.org/obo/ENVO_00002018,https://en.wikipedia.org/wiki/Sewage"
$ws.Range("E180").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13278"
$ws.Range("H180").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13278"
$ws.Range("E181").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C33739"
$ws.Range("H181").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C33739"
$ws.Range("E182").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C12801"
$ws.Range("H182").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C12801"
$ws.Range("E183").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13283"
$ws.Range("H183").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C13283"
$ws.Range("E184").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C173496"
$ws.Range("H184").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C173496"
$ws.Range("E185").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C28294"
$ws.Range("H185").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C28294"
$ws.Range("E186").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C150892"
$ws.Range("H186").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C150892"
$ws.Range("E187").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C19157"
$ws.Range("H187").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C19157,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/123038009"
$ws.Range("B188").Value = "cell line"
$ws.Range("E188").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C16403"
$ws.Range("H188").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C16403"
$ws.Range("E189").Value = "http://purl.obolibrary.org/obo/OBI_0000971"
$ws.Range("E190").Value = "http://purl.obolibrary.org/obo/OBI_0000922"
$ws.Range("D191").Value = "A method used to collect biological material from within the nasal passages. A cotton swab is inserted into the nasal opening and rotated against the anterior nasal mucosa and them withdrawn."
$ws.Range("E191").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C132119"
$ws.Range("H191").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C132119"
$ws.Range("E192").Value = "http://purl.obolibrary.org/obo/IAO_0000122"
$ws.Range("E193").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C155835"
$ws.Range("H193").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C155835,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/461911000124106,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/258529004"
$ws.Range("E194").Value = "http://purl.bioontology.org/ontology/SNOMEDCT/119376003"
$ws.Range("E195").Value = "http://purl.obolibrary.org/obo/OBI_0002611"
$ws.Range("E196").Value = "http://purl.obolibrary.org/obo/OBI_1200000"
$ws.Range("H196").Value = "http://purl.obolibrary.org/obo/OBI_1200000"
$ws.Range("D197").Value = ""
$ws.Range("E197").Value = ""
$ws.Range("H197").Value = "http://purl.obolibrary.org/obo/FBbi_00000011"
$ws.Range("H198").Value = ""
$ws.Range("E199").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C324"
$ws.Range("H199").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C324"
$ws.Range("E200").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C449"
$ws.Range("H200").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C449,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/24851008"
$ws.Range("E201").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C812"
$ws.Range("H201").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C812,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/27888000"
$ws.Range("E202").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C17021"
$ws.Range("H202").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C17021,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/88878007"
$ws.Range("E203").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C16295"
$ws.Range("H203").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C16295,`nhttp://purl.bioontology.org/ontology/SNOMEDCT/68498002"
$ws.Range("E206").Value = "https://www.utwente.nl/.uc/fb9dbdab80102e9e99b00eeb5220284e61d13c337a78900/short-manual-gmo.pdf"
$ws.Range("E207").Value = "https://www.utwente.nl/.uc/fb9dbdab80102e9e99b00eeb5220284e61d13c337a78900/short-manual-gmo.pdf"
$ws.Range("E208").Value = "https://www.utwente.nl/.uc/fb9dbdab80102e9e99b00eeb5220284e61d13c337a78900/short-manual-gmo.pdf"
$ws.Range("E209").Value = "https://www.labmanager.com/lab-health-and-safety/biosafety-levels-1-2-3-4-19123"
$ws.Range("E210").Value = "https://www.labmanager.com/lab-health-and-safety/biosafety-levels-1-2-3-4-19123"
$ws.Range("E211").Value = "https://www.labmanager.com/lab-health-and-safety/biosafety-levels-1-2-3-4-19123"
$ws.Range("E212").Value = "https://www.labmanager.com/lab-health-and-safety/biosafety-levels-1-2-3-4-19123"
$ws.Range("B215").Value = "metabolomics facility"
$ws.Range("B226").Value = "repository"
$ws.Range("F226").Value = "zonmwpc:10186"
$ws.Range("B227").Value = "biobank"
$ws.Range("B228").Value = "data catalogue"
$ws.Range("B229").Value = "ELSI service"
$ws.Range("B230").Value = "other online service"

# Delete row 231 (last row, becomes dimension A1:X230)
$ws.Rows.Item(231).Delete()
